$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix data for "US Core CarePlan Profile" (row 14): the old C14 value
# "Device" (a bare resource type) is replaced by the proper profile
# reference "US Core Device Profile", and the remaining target columns
# shift left by one.
$ws.Cells.Item(14, 3).Value = "PractitionerRole"
$ws.Cells.Item(14, 4).Value = "US Core CareTeam Profile"
$ws.Cells.Item(14, 5).Value = "US Core Device Profile"

# --- Insert a new row for "US Core Device Profile" right after
# "US Core Coverage Profile" (originally row 18), i.e. before the old row 19.
$ws.Rows(19).Insert()
$ws.Cells.Item(19, 2).Value = "US Core Device Profile"

# --- Insert a new row for "US Core FamilyMemberHistory Profile" right
# after "US Core Encounter Profile" (before the old "Goal Profile" row,
# which is now row 24 because of the previous insertion).
$ws.Rows(24).Insert()
$ws.Cells.Item(24, 2).Value = "US Core FamilyMemberHistory Profile"

# --- Remove the row for "US Core Implantable Device Profile" entirely
# (originally row 27, now row 29 after the two insertions above).
$ws.Rows(29).Delete()

# --- Insert a new row for "US Core PMO ServiceRequest Profile" right
# after "US Core Patient Profile" (before the old "Practitioner Profile"
# row, which is now row 44; insert before it at row 43).
$ws.Rows(43).Insert()
$ws.Cells.Item(43, 2).Value = "US Core PMO ServiceRequest Profile"

# --- Renumber column A (the 0-based row index) for every data row so it
# stays sequential after the inserts/deletes above. Column B (the profile
# name) is always populated, so use it to find the last data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
